$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 4 de Agosto de 2020 a las 08:59"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 4862285
$ws.Range("C4").Value = 111
$ws.Range("D4").Value = 2447525
$ws.Range("E4").Value = 2255829
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 158931

# Row 6 - India
$ws.Range("B6").Value = 1858689
$ws.Range("C6").Value = 3358
$ws.Range("D6").Value = 1231682
$ws.Range("E6").Value = 588005
$ws.Range("G6").Value = 31
$ws.Range("H6").Value = 39002

# Row 37 - Ucrania
$ws.Range("B37").Value = 74219
$ws.Range("C37").Value = 1061
$ws.Range("D37").Value = 40613
$ws.Range("E37").Value = 31842
$ws.Range("G37").Value = 26
$ws.Range("H37").Value = 1764

# Row 53 - Armenia
$ws.Range("B53").Value = 39298
$ws.Range("C53").Value = 196
$ws.Range("D53").Value = 30372
$ws.Range("E53").Value = 8158
$ws.Range("G53").Value = 6
$ws.Range("H53").Value = 768

# Row 141 - Letonia
$ws.Range("B141").Value = 1249
$ws.Range("C141").Value = 3
$ws.Range("D141").Value = 1070
$ws.Range("E141").Value = 147

# Row 145 - Georgia
$ws.Range("B145").Value = 1182
$ws.Range("C145").Value = 3
$ws.Range("D145").Value = 962

# Row 165 - Taiwan
$ws.Range("B165").Value = 476
$ws.Range("C165").Value = 1
$ws.Range("E165").Value = 28
